# repo clean and add
# - Insert a new "Thing/Model" row into the Deployment Details sheet for the
#   new r1-test-error-file.xlsx test fixture, renumber the Sequence column,
#   and switch the active sheet/selection to where the edit was made.

$wb = $excel.ActiveWorkbook

$tenant = $wb.Worksheets.Item("Tenant Config")
$deploy = $wb.Worksheets.Item("Deployment Details")

# Insert a new row above the old row 3 (the "r1-test-create-entity.xlsx"
# row), pushing the create/delete-entity rows down by one.
$deploy.Rows("3:3").Insert() | Out-Null

# Give the new row the same look as the other "Thing / Model" rows.
$deploy.Range("A4:C4").Copy() | Out-Null
$deploy.Range("A3:C3").PasteSpecial(-4122) | Out-Null
# The Filename column for this row should not wrap (matches row 2's style).
$deploy.Range("D2").Copy() | Out-Null
$deploy.Range("D3").PasteSpecial(-4122) | Out-Null

# Fill in the new row's content.
$deploy.Range("A3").Value = 1
$deploy.Range("B3").Value = "Thing"
$deploy.Range("C3").Value = "Model"
$deploy.Range("D3").Value = "r1-test-error-file.xlsx"

# Renumber the Sequence column for the rows that shifted / follow it.
$deploy.Range("A2").Value = 2
$deploy.Range("A4").Value = 3
$deploy.Range("A5").Value = 4

# The edits were made on the Deployment Details sheet, so leave it selected.
$deploy.Activate() | Out-Null
$deploy.Range("B10").Select() | Out-Null
$tenant.Range("B19").Select() | Out-Null
